$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Hgf"
$ws.Range("C2").Value = "St14"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 5.231719
$ws.Range("H2").Value = 15.695157
$ws.Range("I2").Value = 0.362499186434781
$ws.Range("J2").Value = 0.362499186434781
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.051909333333333
$ws.Range("N2").Value = 3.155728
$ws.Range("O2").Value = 0.6689916625345733
$ws.Range("P2").Value = 0.6689916625345734
$ws.Range("Q2").Value = 5.503294045477333
$ws.Range("R2").Value = 49.529646409296
$ws.Range("S2").Value = 0.2425089334004344
$ws.Range("T2").Value = 0.2425089334004344

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Hgf"
$ws.Range("C3").Value = "St14"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 5.231719
$ws.Range("H3").Value = 15.695157
$ws.Range("I3").Value = 0.362499186434781
$ws.Range("J3").Value = 0.362499186434781
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1261633333333333
$ws.Range("N3").Value = 0.37849
$ws.Range("O3").Value = 0.08023716060215286
$ws.Range("P3").Value = 0.08023716060215287
$ws.Range("Q3").Value = 0.6600511081033332
$ws.Range("R3").Value = 5.94045997293
$ws.Range("S3").Value = 0.02908590544011727
$ws.Range("T3").Value = 0.02908590544011727

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Hgf"
$ws.Range("C4").Value = "St14"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 5.231719
$ws.Range("H4").Value = 15.695157
$ws.Range("I4").Value = 0.362499186434781
$ws.Range("J4").Value = 0.362499186434781
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3943076666666667
$ws.Range("N4").Value = 1.182923
$ws.Range("O4").Value = 0.2507711768632738
$ws.Range("P4").Value = 0.2507711768632738
$ws.Range("Q4").Value = 2.062906911545667
$ws.Range("R4").Value = 18.566162203911
$ws.Range("S4").Value = 0.09090434759422931
$ws.Range("T4").Value = 0.09090434759422931

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Hgf"
$ws.Range("C5").Value = "St14"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.316945333333335
$ws.Range("H5").Value = 24.950836
$ws.Range("I5").Value = 0.5762706133406404
$ws.Range("J5").Value = 0.5762706133406403
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.051909333333333
$ws.Range("N5").Value = 3.155728
$ws.Range("O5").Value = 0.6689916625345733
$ws.Range("P5").Value = 0.6689916625345734
$ws.Range("Q5").Value = 8.748672420956446
$ws.Range("R5").Value = 78.738051788608
$ws.Range("S5").Value = 0.3855202356885733
$ws.Range("T5").Value = 0.3855202356885733

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hgf"
$ws.Range("C6").Value = "St14"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.316945333333335
$ws.Range("H6").Value = 24.950836
$ws.Range("I6").Value = 0.5762706133406404
$ws.Range("J6").Value = 0.5762706133406403
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1261633333333333
$ws.Range("N6").Value = 0.37849
$ws.Range("O6").Value = 0.08023716060215286
$ws.Range("P6").Value = 0.08023716060215287
$ws.Range("Q6").Value = 1.049293546404445
$ws.Range("R6").Value = 9.44364191764
$ws.Range("S6").Value = 0.0462383177529141
$ws.Range("T6").Value = 0.0462383177529141

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Hgf"
$ws.Range("C7").Value = "St14"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.316945333333335
$ws.Range("H7").Value = 24.950836
$ws.Range("I7").Value = 0.5762706133406404
$ws.Range("J7").Value = 0.5762706133406403
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3943076666666667
$ws.Range("N7").Value = 1.182923
$ws.Range("O7").Value = 0.2507711768632738
$ws.Range("P7").Value = 0.2507711768632738
$ws.Range("Q7").Value = 3.27943530818089
$ws.Range("R7").Value = 29.51491777362801
$ws.Range("S7").Value = 0.144512059899153
$ws.Range("T7").Value = 0.144512059899153

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Hgf"
$ws.Range("C8").Value = "St14"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8836963333333333
$ws.Range("H8").Value = 2.651089
$ws.Range("I8").Value = 0.06123020022457864
$ws.Range("J8").Value = 0.06123020022457864
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.051909333333333
$ws.Range("N8").Value = 3.155728
$ws.Range("O8").Value = 0.6689916625345733
$ws.Range("P8").Value = 0.6689916625345734
$ws.Range("Q8").Value = 0.9295684208657777
$ws.Range("R8").Value = 8.366115787792
$ws.Range("S8").Value = 0.04096249344556567
$ws.Range("T8").Value = 0.04096249344556568

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Hgf"
$ws.Range("C9").Value = "St14"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8836963333333333
$ws.Range("H9").Value = 2.651089
$ws.Range("I9").Value = 0.06123020022457864
$ws.Range("J9").Value = 0.06123020022457864
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1261633333333333
$ws.Range("N9").Value = 0.37849
$ws.Range("O9").Value = 0.08023716060215286
$ws.Range("P9").Value = 0.08023716060215287
$ws.Range("Q9").Value = 0.1114900750677778
$ws.Range("R9").Value = 1.00341067561
$ws.Range("S9").Value = 0.004912937409121492
$ws.Range("T9").Value = 0.004912937409121493

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Hgf"
$ws.Range("C10").Value = "St14"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8836963333333333
$ws.Range("H10").Value = 2.651089
$ws.Range("I10").Value = 0.06123020022457864
$ws.Range("J10").Value = 0.06123020022457864
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3943076666666667
$ws.Range("N10").Value = 1.182923
$ws.Range("O10").Value = 0.2507711768632738
$ws.Range("P10").Value = 0.2507711768632738
$ws.Range("Q10").Value = 0.3484482392385556
$ws.Range("R10").Value = 3.136034153147
$ws.Range("S10").Value = 0.01535476936989148
$ws.Range("T10").Value = 0.01535476936989148
